$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") '40.989.86'
Set-TextValue $ws.Range("E2") '  -1.63%  '
Set-TextValue $ws.Range("D3") '2.164.26'
Set-TextValue $ws.Range("E3") '  -2.62%  '
Set-TextValue $ws.Range("E4") '  -0.02%  '
Set-TextValue $ws.Range("D5") '246.36'
Set-TextValue $ws.Range("E5") '  -1.82%  '
Set-TextValue $ws.Range("D6") '0.611'
Set-TextValue $ws.Range("E6") '  -2.69%  '
Set-TextValue $ws.Range("D7") '65.39'
Set-TextValue $ws.Range("E7") '  -7.82%  '
Set-TextValue $ws.Range("E8") '  +0.07%  '
Set-TextValue $ws.Range("D9") '0.558'
Set-TextValue $ws.Range("E9") '  -3.46%  '
Set-TextValue $ws.Range("D10") '59.50'
Set-TextValue $ws.Range("E10") '  +2.24%  '
Set-TextValue $ws.Range("D11") '0.0918'
Set-TextValue $ws.Range("E11") '  -4.85%  '
Set-TextValue $ws.Range("D12") '34.91'
Set-TextValue $ws.Range("E12") '  -15.58%  '
Set-TextValue $ws.Range("E13") '  -2.01%  '
Set-TextValue $ws.Range("D14") '6.77'
Set-TextValue $ws.Range("E14") '  -3.16%  '
Set-TextValue $ws.Range("D15") '2.489.16'
Set-TextValue $ws.Range("E15") '  -2.49%  '
Set-TextValue $ws.Range("D16") '14.15'
Set-TextValue $ws.Range("E16") '  -4.87%  '
Set-TextValue $ws.Range("D17") '0.842'
Set-TextValue $ws.Range("E17") '  -1.59%  '
Set-TextValue $ws.Range("D18") '2.174.96'
Set-TextValue $ws.Range("E18") '  -2.14%  '
Set-TextValue $ws.Range("D19") '40.873.59'
Set-TextValue $ws.Range("E19") '  -1.67%  '
Set-TextValue $ws.Range("D20") '0.0₃0930'
Set-TextValue $ws.Range("E20") '  -3.80%  '
Set-TextValue $ws.Range("D21") '71.10'
Set-TextValue $ws.Range("E21") '  -1.98%  '
Set-TextValue $ws.Range("D22") '6.02'
Set-TextValue $ws.Range("E22") '  -2.95%  '
Set-TextValue $ws.Range("D23") '227.78'
Set-TextValue $ws.Range("E23") '  -2.73%  '
Set-TextValue $ws.Range("D24") '2.04'
Set-TextValue $ws.Range("E24") '  -7.49%  '
Set-TextValue $ws.Range("E26") '  -5.00%  '
Set-TextValue $ws.Range("D27") '11.12'
Set-TextValue $ws.Range("E27") '  +6.58%  '
Set-TextValue $ws.Range("D28") '2.40'
Set-TextValue $ws.Range("E28") '  -3.98%  '
Set-TextValue $ws.Range("E29") '  -5.61%  '
Set-TextValue $ws.Range("D30") '167.64'
Set-TextValue $ws.Range("E30") '  -2.09%  '
Set-TextValue $ws.Range("E31") '  -8.95%  '
Set-TextValue $ws.Range("D32") '20.05'
Set-TextValue $ws.Range("E32") '  -2.75%  '
Set-TextValue $ws.Range("D33") '0.120'
Set-TextValue $ws.Range("E33") '  -0.42%  '
Set-TextValue $ws.Range("E34") '  +1.45%  '
Set-TextValue $ws.Range("D35") '0.0738'
Set-TextValue $ws.Range("E35") '  +2.63%  '
Set-TextValue $ws.Range("E36") '  -3.52%  '
Set-TextValue $ws.Range("E37") '  -3.19%  '
Set-TextValue $ws.Range("D38") '3.93'
Set-TextValue $ws.Range("E38") '  +0.09%  '
Set-TextValue $ws.Range("D39") '24.18'
Set-TextValue $ws.Range("E39") '  -7.20%  '
Set-TextValue $ws.Range("D40") '0.0298'
Set-TextValue $ws.Range("E40") '  +1.11%  '
Set-TextValue $ws.Range("D41") '2.16'
Set-TextValue $ws.Range("E41") '  -5.40%  '
Set-TextValue $ws.Range("D42") '5.42'
Set-TextValue $ws.Range("E42") '  -8.55%  '
Set-TextValue $ws.Range("D43") '4.81'
Set-TextValue $ws.Range("E43") '  -0.81%  '
Set-TextValue $ws.Range("D44") '59.78'
Set-TextValue $ws.Range("E44") '  -12.75%  '
Set-TextValue $ws.Range("E45") '  -6.68%  '
Set-TextValue $ws.Range("D46") '0.189'
Set-TextValue $ws.Range("E46") '  -9.55%  '
Set-TextValue $ws.Range("E47") '  -0.07%  '
Set-TextValue $ws.Range("D48") '8.41'
Set-TextValue $ws.Range("E48") '  -4.32%  '
Set-TextValue $ws.Range("D49") '0.0981'
Set-TextValue $ws.Range("E49") '  -3.34%  '
Set-TextValue $ws.Range("D50") '1.13'
Set-TextValue $ws.Range("E50") '  -1.27%  '
Set-TextValue $ws.Range("E51") '  -4.40%  '
